$wb = $excel.ActiveWorkbook

# Remove the two extra "delta" sheets
$wb.Worksheets.Item("delta1").Delete() | Out-Null
$wb.Worksheets.Item("delta2").Delete() | Out-Null

# Rename the first sheet "Sheet" -> "Schema"
$schema = $wb.Worksheets.Item("Sheet")
$schema.Name = "Schema"

# Populate the Schema sheet with header + rows describing the delta sheet's schema
$schema.Range("A1").Value = "Sheet"
$schema.Range("B1").Value = "Order"
$schema.Range("C1").Value = "Name"
$schema.Range("D1").Value = "Key"
$schema.Range("E1").Value = "Unique"
$schema.Range("F1").Value = "Foreign Sheet"
$schema.Range("G1").Value = "Foreign Key"

$schema.Range("A2").Value = "delta"
$schema.Range("B2").Value = 0
$schema.Range("C2").Value = "prop_a"

$schema.Range("A3").Value = "delta"
$schema.Range("B3").Value = 1
$schema.Range("C3").Value = "prop_b"

$schema.Range("A4").Value = "delta"
$schema.Range("B4").Value = 2
$schema.Range("C4").Value = "prop_c"

# Update the delta sheet: mark column C (prop_c) as a boolean "unique" flag for each data row
$delta = $wb.Worksheets.Item("delta")
for ($r = 2; $r -le 7; $r++) {
    $delta.Cells.Item($r, 3).Value = $true
}

# Remove the now-excess rows 8 and 9 from delta
$delta.Rows.Item(9).Delete() | Out-Null
$delta.Rows.Item(8).Delete() | Out-Null

# Restore the originally active tab (first sheet)
$schema.Activate()

Write-Output "done"
